$d = $word.ActiveDocument

# The paragraph ending in "LOM3046: Técnicas de Análise Microestrutural
# (Requisito)" is immediately followed by three paragraphs that must be
# removed: a blank paragraph, "Ver no Jupiter Salvar em pdf Salvar em docx",
# and "© 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github
# pages. Original theme under Creative Commons Attribution". Removing them
# leaves the following (already blank) paragraph directly after the LOM3046
# paragraph, right before the page-break paragraph.

$startRng = $d.Content
$startRng.Find.Execute("LOM3046: Técnicas de Análise Microestrutural (Requisito)", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
# Find's Range.End lands right before the paragraph mark of the found
# paragraph; add 1 to move past it, to the start of the next paragraph.
$startPos = $startRng.End + 1

$endRng = $d.Content
$endRng.Find.Execute("© 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
# Likewise, add 1 to include that paragraph's own paragraph mark in the
# deleted range.
$endPos = $endRng.End + 1

$deleteRange = $d.Range($startPos, $endPos)
$deleteRange.Delete()
